# Auto-generated edit script applying the Adamantoise_Profits.xlsx diff
# Each hunk corresponds to specific (sheet, row) leve entries whose
# market-price / profit columns (H:N) were refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 3029.6667  # H51: 3051.3794 -> 3029.6667
$ws.Cells.Item(51, 10).Value = 2938.3333  # J51: 2970 -> 2938.3333
$ws.Cells.Item(51, 12).Value = 2938.3333  # L51: 2970 -> 2938.3333
$ws.Cells.Item(51, 14).Value = -3906.3333  # N51: -3938 -> -3906.3333
$ws.Cells.Item(132, 8).Value = 3634.3542  # H132: 3691.7446 -> 3634.3542
$ws.Cells.Item(132, 9).Value = 3640.9268  # I132: 3708.525 -> 3640.9268
$ws.Cells.Item(132, 11).Value = 10922.7804  # K132: 11125.575 -> 10922.7804
$ws.Cells.Item(132, 13).Value = -8392.7804  # M132: -8595.575000000001 -> -8392.7804
$ws.Cells.Item(138, 8).Value = 2565.0317  # H138: 2576.127 -> 2565.0317
$ws.Cells.Item(138, 9).Value = 2362.4211  # I138: 2449.2222 -> 2362.4211
$ws.Cells.Item(138, 10).Value = 2652.5227  # J138: 2626.889 -> 2652.5227
$ws.Cells.Item(138, 11).Value = 7087.263300000001  # K138: 7347.6666 -> 7087.263300000001
$ws.Cells.Item(138, 12).Value = 7957.5681  # L138: 7880.667 -> 7957.5681
$ws.Cells.Item(138, 13).Value = -1947.263300000001  # M138: -2207.6666 -> -1947.263300000001
$ws.Cells.Item(138, 14).Value = -18237.5681  # N138: -18160.667 -> -18237.5681
$ws.Cells.Item(141, 8).Value = 4236.5713  # H141: 4006.375 -> 4236.5713
$ws.Cells.Item(141, 9).Value = 3626.1667  # I141: 3450.2856 -> 3626.1667
$ws.Cells.Item(141, 11).Value = 10878.5001  # K141: 10350.8568 -> 10878.5001
$ws.Cells.Item(141, 13).Value = -5698.500100000001  # M141: -5170.856800000001 -> -5698.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2715.476  # H45: 4413.273 -> 2715.476
$ws.Cells.Item(45, 9).Value = 1935.8667  # I45: 3221.1428 -> 1935.8667
$ws.Cells.Item(45, 10).Value = 4664.5  # J45: 6499.5 -> 4664.5
$ws.Cells.Item(45, 11).Value = 1935.8667  # K45: 3221.1428 -> 1935.8667
$ws.Cells.Item(45, 12).Value = 4664.5  # L45: 6499.5 -> 4664.5
$ws.Cells.Item(45, 13).Value = -1558.8667  # M45: -2844.1428 -> -1558.8667
$ws.Cells.Item(45, 14).Value = -5418.5  # N45: -7253.5 -> -5418.5
$ws.Cells.Item(61, 8).Value = 2772.6758  # H61: 2817.8333 -> 2772.6758
$ws.Cells.Item(61, 9).Value = 2613.7273  # I61: 2659.5625 -> 2613.7273
$ws.Cells.Item(61, 11).Value = 2613.7273  # K61: 2659.5625 -> 2613.7273
$ws.Cells.Item(61, 13).Value = -2401.7273  # M61: -2447.5625 -> -2401.7273
$ws.Cells.Item(74, 8).Value = 2178  # H74: 2326 -> 2178
$ws.Cells.Item(74, 9).Value = 1330.6842  # I74: 1422.5294 -> 1330.6842
$ws.Cells.Item(74, 11).Value = 1330.6842  # K74: 1422.5294 -> 1330.6842
$ws.Cells.Item(74, 13).Value = -456.6841999999999  # M74: -548.5293999999999 -> -456.6841999999999
$ws.Cells.Item(77, 8).Value = 2178  # H77: 2326 -> 2178
$ws.Cells.Item(77, 9).Value = 1330.6842  # I77: 1422.5294 -> 1330.6842
$ws.Cells.Item(77, 11).Value = 6653.420999999999  # K77: 7112.646999999999 -> 6653.420999999999
$ws.Cells.Item(77, 13).Value = -2285.420999999999  # M77: -2744.646999999999 -> -2285.420999999999
$ws.Cells.Item(122, 8).Value = 4189.0405  # H122: 4317.1406 -> 4189.0405
$ws.Cells.Item(122, 9).Value = 3220.7114  # I122: 3347.0408 -> 3220.7114
$ws.Cells.Item(122, 11).Value = 9662.1342  # K122: 10041.1224 -> 9662.1342
$ws.Cells.Item(122, 13).Value = -7212.1342  # M122: -7591.1224 -> -7212.1342
$ws.Cells.Item(132, 8).Value = 2660.0815  # H132: 2709.4893 -> 2660.0815
$ws.Cells.Item(132, 9).Value = 2206.842  # I132: 2246.1667 -> 2206.842
$ws.Cells.Item(132, 11).Value = 6620.526  # K132: 6738.500100000001 -> 6620.526
$ws.Cells.Item(132, 13).Value = -4090.526  # M132: -4208.500100000001 -> -4090.526
$ws.Cells.Item(136, 8).Value = 2772.6758  # H136: 2817.8333 -> 2772.6758
$ws.Cells.Item(136, 9).Value = 2613.7273  # I136: 2659.5625 -> 2613.7273
$ws.Cells.Item(136, 11).Value = 7841.1819  # K136: 7978.6875 -> 7841.1819
$ws.Cells.Item(136, 13).Value = -5291.1819  # M136: -5428.6875 -> -5291.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3491.7727  # H86: 3884.7368 -> 3491.7727
$ws.Cells.Item(86, 9).Value = 3169  # I86: 3891 -> 3169
$ws.Cells.Item(86, 11).Value = 3169  # K86: 3891 -> 3169
$ws.Cells.Item(86, 13).Value = -2046  # M86: -2768 -> -2046
$ws.Cells.Item(89, 8).Value = 3491.7727  # H89: 3884.7368 -> 3491.7727
$ws.Cells.Item(89, 9).Value = 3169  # I89: 3891 -> 3169
$ws.Cells.Item(89, 11).Value = 15845  # K89: 19455 -> 15845
$ws.Cells.Item(89, 13).Value = -10229  # M89: -13839 -> -10229
$ws.Cells.Item(134, 8).Value = 13495513  # H134: 13960859 -> 13495513
$ws.Cells.Item(134, 9).Value = 2859282.2  # I134: 2978399 -> 2859282.2
$ws.Cells.Item(134, 11).Value = 8577846.600000001  # K134: 8935197 -> 8577846.600000001
$ws.Cells.Item(134, 13).Value = -8575311.600000001  # M134: -8932662 -> -8575311.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3755.9048  # H31: 4023.15 -> 3755.9048
$ws.Cells.Item(31, 9).Value = 2325.5833  # I31: 2408.7273 -> 2325.5833
$ws.Cells.Item(31, 10).Value = 5663  # J31: 5996.3335 -> 5663
$ws.Cells.Item(31, 11).Value = 2325.5833  # K31: 2408.7273 -> 2325.5833
$ws.Cells.Item(31, 12).Value = 5663  # L31: 5996.3335 -> 5663
$ws.Cells.Item(31, 13).Value = -2030.5833  # M31: -2113.7273 -> -2030.5833
$ws.Cells.Item(31, 14).Value = -6253  # N31: -6586.3335 -> -6253
$ws.Cells.Item(34, 8).Value = 3755.9048  # H34: 4023.15 -> 3755.9048
$ws.Cells.Item(34, 9).Value = 2325.5833  # I34: 2408.7273 -> 2325.5833
$ws.Cells.Item(34, 10).Value = 5663  # J34: 5996.3335 -> 5663
$ws.Cells.Item(34, 11).Value = 2325.5833  # K34: 2408.7273 -> 2325.5833
$ws.Cells.Item(34, 12).Value = 5663  # L34: 5996.3335 -> 5663
$ws.Cells.Item(34, 13).Value = -2123.5833  # M34: -2206.7273 -> -2123.5833
$ws.Cells.Item(34, 14).Value = -6067  # N34: -6400.3335 -> -6067
$ws.Cells.Item(99, 8).Value = 3252.75  # H99: 3374.75 -> 3252.75
$ws.Cells.Item(99, 9).Value = 1837  # I99: 1999.6666 -> 1837
$ws.Cells.Item(99, 11).Value = 1837  # K99: 1999.6666 -> 1837
$ws.Cells.Item(99, 13).Value = -339  # M99: -501.6666 -> -339
$ws.Cells.Item(105, 8).Value = 2731.2856  # H105: 2895.8 -> 2731.2856
$ws.Cells.Item(105, 9).Value = 2323.8  # I105: 2326.3333 -> 2323.8
$ws.Cells.Item(105, 11).Value = 2323.8  # K105: 2326.3333 -> 2323.8
$ws.Cells.Item(105, 13).Value = -576.8000000000002  # M105: -579.3332999999998 -> -576.8000000000002
$ws.Cells.Item(126, 8).Value = 3252.75  # H126: 3374.75 -> 3252.75
$ws.Cells.Item(126, 9).Value = 1837  # I126: 1999.6666 -> 1837
$ws.Cells.Item(126, 11).Value = 5511  # K126: 5998.9998 -> 5511
$ws.Cells.Item(126, 13).Value = -3041  # M126: -3528.9998 -> -3041
$ws.Cells.Item(132, 8).Value = 1856.9032  # H132: 1748.8 -> 1856.9032
$ws.Cells.Item(132, 9).Value = 1752.4615  # I132: 1665.8276 -> 1752.4615
$ws.Cells.Item(132, 10).Value = 2400  # J132: 2149.8333 -> 2400
$ws.Cells.Item(132, 11).Value = 5257.3845  # K132: 4997.4828 -> 5257.3845
$ws.Cells.Item(132, 12).Value = 7200  # L132: 6449.499899999999 -> 7200
$ws.Cells.Item(132, 13).Value = -2727.3845  # M132: -2467.4828 -> -2727.3845
$ws.Cells.Item(132, 14).Value = -12260  # N132: -11509.4999 -> -12260
$ws.Cells.Item(134, 8).Value = 2940.25  # H134: 2870.476 -> 2940.25
$ws.Cells.Item(134, 9).Value = 2147.7693  # I134: 2099.7144 -> 2147.7693
$ws.Cells.Item(134, 11).Value = 6443.3079  # K134: 6299.1432 -> 6443.3079
$ws.Cells.Item(134, 13).Value = -3908.3079  # M134: -3764.1432 -> -3908.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 17  # H26: 17.75 -> 17
$ws.Cells.Item(26, 10).Value = 1  # J26: 4 -> 1
$ws.Cells.Item(26, 12).Value = 3  # L26: 12 -> 3
$ws.Cells.Item(26, 14).Value = -579  # N26: -588 -> -579

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2777.4255  # H132: 3150.6584 -> 2777.4255
$ws.Cells.Item(132, 9).Value = 1887.303  # I132: 2141.4443 -> 1887.303
$ws.Cells.Item(132, 10).Value = 4875.5713  # J132: 5097 -> 4875.5713
$ws.Cells.Item(132, 11).Value = 5661.909000000001  # K132: 6424.3329 -> 5661.909000000001
$ws.Cells.Item(132, 12).Value = 14626.7139  # L132: 15291 -> 14626.7139
$ws.Cells.Item(132, 13).Value = -3131.909000000001  # M132: -3894.3329 -> -3131.909000000001
$ws.Cells.Item(132, 14).Value = -19686.7139  # N132: -20351 -> -19686.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4440.421  # H7: 4549.8887 -> 4440.421
$ws.Cells.Item(7, 9).Value = 3668.1  # I7: 3711.1 -> 3668.1
$ws.Cells.Item(7, 10).Value = 5298.5557  # J7: 5598.375 -> 5298.5557
$ws.Cells.Item(7, 11).Value = 3668.1  # K7: 3711.1 -> 3668.1
$ws.Cells.Item(7, 12).Value = 5298.5557  # L7: 5598.375 -> 5298.5557
$ws.Cells.Item(7, 13).Value = -3556.1  # M7: -3599.1 -> -3556.1
$ws.Cells.Item(7, 14).Value = -5522.5557  # N7: -5822.375 -> -5522.5557
$ws.Cells.Item(40, 8).Value = 6151.6924  # H40: 6397.923 -> 6151.6924
$ws.Cells.Item(40, 9).Value = 5965.84  # I40: 6221.92 -> 5965.84
$ws.Cells.Item(40, 11).Value = 5965.84  # K40: 6221.92 -> 5965.84
$ws.Cells.Item(40, 13).Value = -5829.84  # M40: -6085.92 -> -5829.84
$ws.Cells.Item(126, 8).Value = 4440.421  # H126: 4549.8887 -> 4440.421
$ws.Cells.Item(126, 9).Value = 3668.1  # I126: 3711.1 -> 3668.1
$ws.Cells.Item(126, 10).Value = 5298.5557  # J126: 5598.375 -> 5298.5557
$ws.Cells.Item(126, 11).Value = 11004.3  # K126: 11133.3 -> 11004.3
$ws.Cells.Item(126, 12).Value = 15895.6671  # L126: 16795.125 -> 15895.6671
$ws.Cells.Item(126, 13).Value = -8534.299999999999  # M126: -8663.299999999999 -> -8534.299999999999
$ws.Cells.Item(126, 14).Value = -20835.6671  # N126: -21735.125 -> -20835.6671
$ws.Cells.Item(132, 8).Value = 2620.818  # H132: 2418.1072 -> 2620.818
$ws.Cells.Item(132, 9).Value = 2292.5  # I132: 2148.5217 -> 2292.5
$ws.Cells.Item(132, 10).Value = 4098.25  # J132: 3658.2 -> 4098.25
$ws.Cells.Item(132, 11).Value = 6877.5  # K132: 6445.5651 -> 6877.5
$ws.Cells.Item(132, 12).Value = 12294.75  # L132: 10974.6 -> 12294.75
$ws.Cells.Item(132, 13).Value = -4347.5  # M132: -3915.5651 -> -4347.5
$ws.Cells.Item(132, 14).Value = -17354.75  # N132: -16034.6 -> -17354.75
$ws.Cells.Item(136, 8).Value = 3327.7908  # H136: 3150.0625 -> 3327.7908
$ws.Cells.Item(136, 9).Value = 2857.6333  # I136: 2681.0571 -> 2857.6333
$ws.Cells.Item(136, 11).Value = 8572.8999  # K136: 8043.1713 -> 8572.8999
$ws.Cells.Item(136, 13).Value = -6022.8999  # M136: -5493.1713 -> -6022.8999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 2000  # H17: 0 -> 2000
$ws.Cells.Item(17, 9).Value = 2000  # I17: 0 -> 2000
$ws.Cells.Item(17, 11).Value = 2000  # K17: 0 -> 2000
$ws.Cells.Item(17, 13).Value = -1828  # M17: None -> -1828
$ws.Cells.Item(19, 8).Value = 18000  # H19: 10005 -> 18000
$ws.Cells.Item(19, 9).Value = 18000  # I19: 10005 -> 18000
$ws.Cells.Item(19, 11).Value = 18000  # K19: 10005 -> 18000
$ws.Cells.Item(19, 13).Value = -17826  # M19: -9831 -> -17826
$ws.Cells.Item(20, 8).Value = 0  # H20: 22010.5 -> 0
$ws.Cells.Item(20, 9).Value = 0  # I20: 14010 -> 0
$ws.Cells.Item(20, 10).Value = 0  # J20: 30011 -> 0
$ws.Cells.Item(20, 11).Value = 0  # K20: 14010 -> 0
$ws.Cells.Item(20, 12).Value = $null  # L20: 30011 -> None
$ws.Cells.Item(20, 13).Value = $null  # M20: -13770 -> None
$ws.Cells.Item(20, 14).Value = 0  # N20: -30491 -> 0
$ws.Cells.Item(23, 8).Value = 999  # H23: 2332.6667 -> 999
$ws.Cells.Item(23, 9).Value = 999  # I23: 999.5 -> 999
$ws.Cells.Item(23, 10).Value = 0  # J23: 4999 -> 0
$ws.Cells.Item(23, 11).Value = 999  # K23: 999.5 -> 999
$ws.Cells.Item(23, 12).Value = 0  # L23: 4999 -> 0
$ws.Cells.Item(23, 13).Value = $null  # M23: -770.5 -> None
$ws.Cells.Item(23, 14).Value = -770  # N23: -5457 -> -770
$ws.Cells.Item(25, 8).Value = 30007.715  # H25: 0 -> 30007.715
$ws.Cells.Item(25, 9).Value = 31250  # I25: 0 -> 31250
$ws.Cells.Item(25, 10).Value = 28351.334  # J25: 0 -> 28351.334
$ws.Cells.Item(25, 11).Value = 31250  # K25: 0 -> 31250
$ws.Cells.Item(25, 12).Value = 28351.334  # L25: 0 -> 28351.334
$ws.Cells.Item(25, 13).Value = -30957  # M25: None -> -30957
$ws.Cells.Item(25, 14).Value = -28937.334  # N25: None -> -28937.334
$ws.Cells.Item(26, 8).Value = 5000  # H26: 7500 -> 5000
$ws.Cells.Item(26, 9).Value = 5000  # I26: 4833.3335 -> 5000
$ws.Cells.Item(26, 10).Value = 0  # J26: 10166.667 -> 0
$ws.Cells.Item(26, 11).Value = 5000  # K26: 4833.3335 -> 5000
$ws.Cells.Item(26, 12).Value = 0  # L26: 10166.667 -> 0
$ws.Cells.Item(26, 13).Value = $null  # M26: -4540.3335 -> None
$ws.Cells.Item(26, 14).Value = -4707  # N26: -10752.667 -> -4707
$ws.Cells.Item(29, 8).Value = 0  # H29: 9999 -> 0
$ws.Cells.Item(29, 10).Value = 0  # J29: 9999 -> 0
$ws.Cells.Item(29, 12).Value = $null  # L29: 9999 -> None
$ws.Cells.Item(29, 14).Value = 0  # N29: -10579 -> 0
$ws.Cells.Item(39, 8).Value = 32596.666  # H39: 30708.5 -> 32596.666
$ws.Cells.Item(39, 9).Value = 32596.666  # I39: 30708.5 -> 32596.666
$ws.Cells.Item(39, 11).Value = 32596.666  # K39: 30708.5 -> 32596.666
$ws.Cells.Item(39, 13).Value = -32183.666  # M39: -30295.5 -> -32183.666
$ws.Cells.Item(54, 8).Value = 48348.332  # H54: 25391.666 -> 48348.332
$ws.Cells.Item(54, 9).Value = 23820  # I54: 22670 -> 23820
$ws.Cells.Item(54, 10).Value = 60612.5  # J54: 39000 -> 60612.5
$ws.Cells.Item(54, 11).Value = 23820  # K54: 22670 -> 23820
$ws.Cells.Item(54, 12).Value = 60612.5  # L54: 39000 -> 60612.5
$ws.Cells.Item(54, 13).Value = -23300  # M54: -22150 -> -23300
$ws.Cells.Item(54, 14).Value = -61652.5  # N54: -40040 -> -61652.5
$ws.Cells.Item(81, 8).Value = 4642.3  # H81: 5561.125 -> 4642.3
$ws.Cells.Item(81, 9).Value = 2346.1428  # I81: 2897.8 -> 2346.1428
$ws.Cells.Item(81, 11).Value = 4692.2856  # K81: 5795.6 -> 4692.2856
$ws.Cells.Item(81, 13).Value = -3631.2856  # M81: -4734.6 -> -3631.2856
$ws.Cells.Item(84, 8).Value = 4642.3  # H84: 5561.125 -> 4642.3
$ws.Cells.Item(84, 9).Value = 2346.1428  # I84: 2897.8 -> 2346.1428
$ws.Cells.Item(84, 11).Value = 23461.428  # K84: 28978 -> 23461.428
$ws.Cells.Item(84, 13).Value = -18157.428  # M84: -23674 -> -18157.428
$ws.Cells.Item(107, 8).Value = 575  # H107: 620.8333 -> 575
$ws.Cells.Item(113, 8).Value = 473.33334  # H113: 473.75 -> 473.33334
$ws.Cells.Item(113, 10).Value = 366.42856  # J113: 367.14285 -> 366.42856
$ws.Cells.Item(113, 12).Value = 1099.28568  # L113: 1101.42855 -> 1099.28568
$ws.Cells.Item(113, 14).Value = -5439.28568  # N113: -5441.428550000001 -> -5439.28568
$ws.Cells.Item(132, 8).Value = 2212.6924  # H132: 2123.9636 -> 2212.6924
$ws.Cells.Item(132, 9).Value = 1638.091  # I132: 1668.721 -> 1638.091
$ws.Cells.Item(132, 10).Value = 5373  # J132: 3755.25 -> 5373
$ws.Cells.Item(132, 11).Value = 4914.272999999999  # K132: 5006.163 -> 4914.272999999999
$ws.Cells.Item(132, 12).Value = 16119  # L132: 11265.75 -> 16119
$ws.Cells.Item(132, 13).Value = -2384.272999999999  # M132: -2476.163 -> -2384.272999999999
$ws.Cells.Item(132, 14).Value = -21179  # N132: -16325.75 -> -21179
$ws.Cells.Item(136, 8).Value = 2810  # H136: 2825.4082 -> 2810
$ws.Cells.Item(136, 9).Value = 1910.7059  # I136: 1961.2122 -> 1910.7059
$ws.Cells.Item(136, 10).Value = 4848.4  # J136: 4607.8125 -> 4848.4
$ws.Cells.Item(136, 11).Value = 5732.1177  # K136: 5883.6366 -> 5732.1177
$ws.Cells.Item(136, 12).Value = 14545.2  # L136: 13823.4375 -> 14545.2
$ws.Cells.Item(136, 13).Value = -3182.1177  # M136: -3333.6366 -> -3182.1177
$ws.Cells.Item(136, 14).Value = -19645.2  # N136: -18923.4375 -> -19645.2

Write-Host "Applied 217 cell updates across 8 sheets"